$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.773749232292175
$ws.Range("B1").Value = 2.538293838500977
$ws.Range("C1").Value = 4.78486967086792
$ws.Range("D1").Value = 4.142022132873535
$ws.Range("E1").Value = 0.9753129482269287
